$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-22 Monday" "2024-07-23 Tuesday"

Replace-Text "94×36=3384" "28×46=1288"
Replace-Text "29×17=493" "65×61=3965"
Replace-Text "85×32=2720" "99×98=9702"
Replace-Text "26×93=2418" "82×33=2706"
Replace-Text "28×91=2548" "53×27=1431"

Replace-Text "16×98=1568" "29×40=1160"
Replace-Text "62×17=1054" "81×61=4941"
Replace-Text "16×35=560" "44×30=1320"
Replace-Text "91×66=6006" "51×19=969"
Replace-Text "31×92=2852" "44×72=3168"

Replace-Text "84×87=7308" "56×62=3472"
Replace-Text "33×25=825" "37×93=3441"
Replace-Text "80×33=2640" "73×65=4745"
Replace-Text "53×62=3286" "65×94=6110"
Replace-Text "25×79=1975" "84×23=1932"

Replace-Text "85×59=5015" "65×77=5005"
Replace-Text "34×92=3128" "74×21=1554"
Replace-Text "38×57=2166" "30×33=990"
Replace-Text "31×18=558" "78×56=4368"
Replace-Text "72×73=5256" "21×28=588"

Replace-Text "74×70=5180" "24×42=1008"
Replace-Text "57×18=1026" "15×80=1200"
Replace-Text "86×22=1892" "73×11=803"
Replace-Text "44×40=1760" "73×24=1752"
Replace-Text "17×45=765" "51×49=2499"
